$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.553958773612976
$ws.Range("B1").Value = 2.737058401107788
$ws.Range("C1").Value = 3.142816305160522
$ws.Range("D1").Value = 2.933822393417358
$ws.Range("E1").Value = 2.761942863464355
